# Historico.xlsx update
# - Insert a new "Id" column before "Area"
# - Replace the table contents with the refreshed dataset (5 data rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at A (shifts Area..Equipo from A:H to B:I)
$ws.Columns.Item(1).Insert()

# 2. Drop the last data row (table now has 5 data rows instead of 6)
$ws.Rows.Item(7).Delete()

# 3. Copy formatting (font/border/fill) from the neighbouring column into
#    the new column A so header + data cells keep the same look.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Match column A's width to the rest of the table (raw width 30).
$ws.Columns.Item(1).ColumnWidth = 29.17

# 5. Header row
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Area"
$ws.Range("C1").Value = "Categoria"
$ws.Range("D1").Value = "Proyecto"
$ws.Range("E1").Value = "Calificación"
$ws.Range("F1").Value = "Lider"
$ws.Range("G1").Value = "Video"
$ws.Range("H1").Value = "Poster"
$ws.Range("I1").Value = "Equipo"

# 6. Data rows (Id, Area, Categoria, Proyecto, Calificación, Lider, Video, Poster, Equipo)
$rows = @(
  @("BioPro1", "Bio",   "Prototipo", "Titulo para mi proyecto de prubas p",      "", "Gerardo Deustúa Hernández", "sdf",                                        "sdf",                                        "Titulo para mi proyecto de prubas p"),
  @("NexPro2", "Nexus", "Producto",  "sfsdfbdsfb",                               "", "Gerardo Deustúa Hernández", "Robot automata para automatizar automatas", "Robot automata para automatizar automatas", "sfsdfbdsfb"),
  @("NexCon1", "Nexus", "Concepto",  "Robot automata para automatizar automatas","", "Gerardo Deustúa Hernández", "sdgasdgasdg",                                "sadgsadg",                                   "Robot automata para automatizar automatas"),
  @("NanPro2", "Nano",  "Prototipo", "Titulo para mi proyecto de prubas Gerry",  "", "Gerardo Deustúa Hernández", "BOLDBGOSDBGOSANGSIGNSOPGSDG",               "BOLDBGOSDBGOSANGSIGNSOPGSDG",               "Titulo para mi proyecto de prubas Gerry"),
  @("NanCon1", "Nano",  "Concepto",  "Test89",                                   "", "Mikel Edel",                "rrrrrrrrrr",                                 "aaaaaaa",                                    "Test89")
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
  $ws.Cells.Item($r, 9).Value = $row[8]
  $r = $r + 1
}

Write-Output "Historico table refreshed"
